# Auto-generated edit script: updates "想去人数" (F column) wishlist counts
# across sheets, plus two cells becoming "暂时售罄" (temporarily sold out) text.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 922
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = 477
$ws.Range("F7").Value = 77
$ws.Range("F8").Value = 0
$ws.Range("F11").Value = 146
$ws.Range("F13").Value = 665
$ws.Range("F14").Value = 529
$ws.Range("F15").Value = 92
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = 34
$ws.Range("F22").Value = 192
$ws.Range("F23").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F29").Value = 31
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 354
$ws.Range("F33").Value = 800
$ws.Range("F34").Value = 348
$ws.Range("F37").Value = 928
$ws.Range("F38").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 317
$ws.Range("F42").Value = 0
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F11").Value = 81
$ws.Range("F12").Value = 6
$ws.Range("F15").Value = 158
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1694
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1694
$ws.Range("F4").Value = 203
$ws.Range("F6").Value = 922
$ws.Range("F7").Value = 83
$ws.Range("F9").Value = 3
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 321
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 8064
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 529
$ws.Range("F22").Value = 81
$ws.Range("F24").Value = 177
$ws.Range("F26").Value = 0
$ws.Range("F28").Value = 525
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 368
$ws.Range("F33").Value = 514
$ws.Range("F34").Value = 0
$ws.Range("F37").Value = 354
$ws.Range("F39").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 171
$ws.Range("F45").Value = 0

# G-column cells that flip from a numeric min price to a "temporarily sold out" label
$ws = $wb.Worksheets.Item("展览")
$ws.Range("G9").Value = "暂时售罄"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G11").Value = "暂时售罄"

